$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.237.32'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.29'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').Value = '  +0.66%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.47'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.008'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5132'
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('E8').Value = '  +1.78%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08382'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.123'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.68'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.247'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.883.40'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.60'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.269'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.010'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001104'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '91.08'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06696'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.79'
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.009'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.027'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.270.99'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.15'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.257'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.097.42'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '160.09'
$ws.Range('E27').Value = '  +1.45%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.473'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '20.69'
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.58'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1056'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.043'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.867'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.612'
$ws.Range('E34').Value = '  +0.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.501'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02445'
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06574'
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2215'
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.199'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6489'
$ws.Range('E41').Value = '  +2.37%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.000'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.23'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6097'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.07'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.697'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.283'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.018'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.238'
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '121.05'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06923'
$ws.Range('E51').Value = '  +1.21%  '
